# Generate Report for Handoff
# Updates Priority ("low" -> "ht") and Latest Handoff Datetime for the rows
# that were just handed off, on both the zh-cn and de-de status sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$zhcnHandoffTime = "2016-08-28 08:31:03"
$dedeHandoffTime = "2016-08-28 08:31:10"

foreach ($row in 4..7) {
    $zhcn.Range("E$row").Value = "ht"
    $zhcn.Range("H$row").Value = $zhcnHandoffTime

    $dede.Range("E$row").Value = "ht"
    $dede.Range("H$row").Value = $dedeHandoffTime
}
